# Insert a new weekly price-report row at row 314 (Ajo / Chino / Primera, Agro
# Chillán), pushing the existing rows 314-380 down to 315-381.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(314).Insert()

$ws.Cells.Item(314, 1).Value = 7
$ws.Cells.Item(314, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(314, 3).Value = "Ñuble"
$ws.Cells.Item(314, 4).Value = 45015
$ws.Cells.Item(314, 5).Value = 16
$ws.Cells.Item(314, 6).Value = 100112003
$ws.Cells.Item(314, 7).Value = "Ajo"
$ws.Cells.Item(314, 8).Value = "Chino"
$ws.Cells.Item(314, 9).Value = "Primera"
$ws.Cells.Item(314, 10).Value = 60
$ws.Cells.Item(314, 11).Value = 17000
$ws.Cells.Item(314, 12).Value = 18000
$ws.Cells.Item(314, 13).Value = 17500
$ws.Cells.Item(314, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(314, 15).Value = "China"
$ws.Cells.Item(314, 16).Value = 1750
$ws.Cells.Item(314, 17).Value = 10
$ws.Cells.Item(314, 18).Value = "Hortaliza"
